# Add a new "2020" column (Q) to the GDP electric-intensity table, mirroring
# the formatting of the existing 2019 column (P).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 gets a taller, explicit height in the new layout.
$ws.Rows.Item(1).RowHeight = 19.5

# --- Copy the formatting from column P into column Q for every row that ---
# --- carries a value/style in P (rows 3 through 8). ------------------------
$ws.Range("P3:P8").Copy() | Out-Null
$ws.Range("Q3:Q8").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- Populate the new 2020 data points. ---
$ws.Range("Q4").Value = 2020
$ws.Range("Q5").Value = 25.6
$ws.Range("Q6").Value = 13.073527219449954
$ws.Range("Q7").Value = 21.941290626870046
$ws.Range("Q8").Value = 196.6
